# Updates the crypto price/volume table to match the latest scrape.
# (Updated cryptos list via GitHub Actions.)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($Cell, $Text) {
    # Round-trip through a Text number format so Excel keeps the literal
    # string (e.g. "5.38") instead of silently coercing it to a number,
    # then clear the format again so no stray style sticks to the cell.
    $ws.Range($Cell).NumberFormat = "@"
    $ws.Range($Cell).Value = $Text
    $ws.Range($Cell).ClearFormats()
}

# Row 2
$ws.Range("D2").Value = "37.819.50"
$ws.Range("E2").Value = "  -0.06%  "

# Row 3
$ws.Range("D3").Value = "2.080.14"

# Row 4
$ws.Range("E4").Value = "  +0.06%  "

# Row 5
Set-TextValue "D5" "233.60"
$ws.Range("E5").Value = "  -0.27%  "

# Row 6
Set-TextValue "D6" "0.624"
$ws.Range("E6").Value = "  -0.10%  "

# Row 7
Set-TextValue "D7" "58.49"
$ws.Range("E7").Value = "  -0.67%  "

# Row 8
$ws.Range("E8").Value = "  -0.01%  "

# Row 9
Set-TextValue "D9" "0.394"
$ws.Range("E9").Value = "  +0.56%  "

# Row 10
Set-TextValue "D10" "0.0785"
$ws.Range("E10").Value = "  -0.69%  "

# Row 11
$ws.Range("E11").Value = "  +3.14%  "

# Row 12
$ws.Range("E12").Value = "  +1.55%  "

# Row 13
$ws.Range("D13").Value = "2.387.10"

# Row 14
Set-TextValue "D14" "21.28"
$ws.Range("E14").Value = "  +0.24%  "

# Row 15
$ws.Range("E15").Value = "  +1.53%  "

# Row 16
Set-TextValue "D16" "5.38"
$ws.Range("E16").Value = "  +1.51%  "

# Row 17
$ws.Range("D17").Value = "2.073.11"
$ws.Range("E17").Value = "  -0.44%  "

# Row 18
$ws.Range("D18").Value = "37.795.41"
$ws.Range("E18").Value = "  +0.10%  "

# Row 19
Set-TextValue "D19" "6.12"
$ws.Range("E19").Value = "  -1.08%  "

# Row 20
Set-TextValue "D20" "71.31"
$ws.Range("E20").Value = "  -0.03%  "

# Row 21
$ws.Range("D21").Value = "0.0₃0838"
$ws.Range("E21").Value = "  +0.47%  "

# Row 22
Set-TextValue "D22" "230.18"
$ws.Range("E22").Value = "  +0.52%  "

# Row 23
$ws.Range("E23").Value = "  -0.12%  "

# Row 24
$ws.Range("E24").Value = "  -0.97%  "

# Row 25
$ws.Range("E25").Value = "  +0.99%  "

# Row 26
Set-TextValue "D26" "9.83"
$ws.Range("E26").Value = "  +9.04%  "

# Row 27
Set-TextValue "D27" "171.93"
$ws.Range("E27").Value = "  +0.84%  "

# Row 28
$ws.Range("E28").Value = "  -2.66%  "

# Row 29
$ws.Range("B29").Value = "ImmutableX"
$ws.Range("C29").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue "D29" "1.41"
$ws.Range("E29").Value = "  -0.26%  "

# Row 30
$ws.Range("B30").Value = "EthereumClassic"
$ws.Range("C30").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextValue "D30" "19.48"
$ws.Range("E30").Value = "  -0.34%  "

# Row 31
$ws.Range("E31").Value = "  +1.20%  "

# Row 32
$ws.Range("E32").Value = "  +0.45%  "

# Row 33
$ws.Range("E33").Value = "  +0.69%  "

# Row 34
$ws.Range("E34").Value = "  -0.29%  "

# Row 35
Set-TextValue "D35" "2.47"
$ws.Range("E35").Value = "  -1.53%  "

# Row 36
$ws.Range("E36").Value = "  -0.57%  "

# Row 37
$ws.Range("E37").Value = "  -2.00%  "

# Row 38
$ws.Range("E38").Value = "  +0.03%  "

# Row 39
$ws.Range("E39").Value = "  +0.92%  "

# Row 40
$ws.Range("E40").Value = "  +9.94%  "

# Row 41
Set-TextValue "D41" "102.26"
$ws.Range("E41").Value = "  +3.55%  "

# Row 42
Set-TextValue "D42" "0.0972"
$ws.Range("E42").Value = "  -1.78%  "

# Row 43
$ws.Range("E43").Value = "  -0.87%  "

# Row 44
$ws.Range("E44").Value = "  +4.34%  "

# Row 45
$ws.Range("D45").Value = "1.450.78"
$ws.Range("E45").Value = "  -1.10%  "

# Row 46
Set-TextValue "D46" "1.15"
$ws.Range("E46").Value = "  -1.20%  "

# Row 47
Set-TextValue "D47" "1.06"
$ws.Range("E47").Value = "  -0.70%  "

# Row 48
$ws.Range("E48").Value = "  -7.14%  "

# Row 49
$ws.Range("E49").Value = "  -0.80%  "

# Row 50
Set-TextValue "D50" "2.99"
$ws.Range("E50").Value = "  -1.33%  "

# Row 51
$ws.Range("D51").Value = "2.271.61"
$ws.Range("E51").Value = "  -0.21%  "
